$d = $word.ActiveDocument

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*rounds*" -and $t -like "*delta*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "target paragraph not found"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi/><w:ind w:left="360"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/></w:rPr><w:lastRenderedPageBreak/><w:t>מכאן התוכנית פולטת</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:ind w:left="360"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>The encryption params are 8 (rounds) and 117132623 (delta)</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:ind w:left="360"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl/></w:rPr><w:t>ומחזירה 0</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:bidi/><w:ind w:left="360"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p = $d.Paragraphs($targetIndex)
$rng = $p.Range
[void]$rng.InsertXML($xml)

Write-Output "done, targetIndex=$targetIndex, newParaCount=$($d.Paragraphs.Count)"
